$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

# Add the new test row (row 18) with the two new values.
$ws.Range("A18").Value = "testRowValue"
$ws.Range("B18").Value = "testRowData."

# Move the active selection to the newly added cell, matching the saved view state.
$ws.Range("B18").Select() | Out-Null
